$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# Rename the "CaseStudy" column header to "CaseStudyName"
$ws.Range("A1").Value = "CaseStudyName"

# Narrow column A to make room for the new municipal water providers layer info
$ws.Columns.Item(1).ColumnWidth = 27.83

# Update the active selection on the Data sheet
$ws.Activate()
$ws.Range("A3").Select()
